# Update the "想去人数" (interest count) figures on both the "展览" and
# "全部类型" worksheets to reflect the newly generated output.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 148
    $ws.Range("F3").Value = 36
    $ws.Range("F4").Value = 230
    $ws.Range("F5").Value = 3833
}
